# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''55.694.01'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.80%  '
$ws.Range("D3").Value = '''2.912.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.33%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '''499.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("D6").Value = '''132.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.89%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '''0.421'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("D9").Value = '''7.13'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.29%  '
$ws.Range("E10").Value = '  -5.50%  '
$ws.Range("E11").Value = '  -4.53%  '
$ws.Range("D12").Value = '''3.410.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.77%  '
$ws.Range("E13").Value = '  -4.08%  '
$ws.Range("D14").Value = '''25.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.42%  '
$ws.Range("E15").Value = '  -3.86%  '
$ws.Range("D16").Value = '''55.609.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.18%  '
$ws.Range("E17").Value = '  -4.67%  '
$ws.Range("D18").Value = '''2.910.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.68%  '
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("D20").Value = '''7.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.97%  '
$ws.Range("D21").Value = '''312.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.21%  '
$ws.Range("D22").Value = '''1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  -2.63%  '
$ws.Range("D24").Value = '''62.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.67%  '
$ws.Range("D25").Value = '''3.029.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.86%  '
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("E27").Value = '  -4.83%  '
$ws.Range("D28").Value = '''0.0₃0825'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -10.50%  '
$ws.Range("D29").Value = '''6.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.42%  '
$ws.Range("E30").Value = '  -9.76%  '
$ws.Range("E31").Value = '  -3.71%  '
$ws.Range("D32").Value = '''19.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.07%  '
$ws.Range("E33").Value = '  -5.92%  '
$ws.Range("D34").Value = '''151.25'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.29%  '
$ws.Range("E35").Value = '  -7.89%  '
$ws.Range("E36").Value = '  -5.50%  '
$ws.Range("D37").Value = '''23.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("D38").Value = '''1.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.42%  '
$ws.Range("D39").Value = '''0.0640'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.99%  '
$ws.Range("D40").Value = '''36.36'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.66%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("D42").Value = '''3.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.31%  '
$ws.Range("D43").Value = '''0.636'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.49%  '
$ws.Range("D44").Value = '''5.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("D45").Value = '''2.102.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.63%  '
$ws.Range("E46").Value = '  -6.15%  '
$ws.Range("D47").Value = '''0.915'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.01%  '
$ws.Range("E48").Value = '  -3.04%  '
$ws.Range("D49").Value = '''18.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.27%  '
$ws.Range("D50").Value = '''0.0837'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.11%  '
$ws.Range("E51").Value = '  -11.05%  '
